$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 125.92857
$ws.Range("I2").Value = 140.11111
$ws.Range("J2").Value = 100.4
$ws.Range("K2").Value = 140.11111
$ws.Range("L2").Value = 100.4
$ws.Range("M2").Value = -27.11111
$ws.Range("N2").Value = -326.4

$ws.Range("H12").Value = 416.7
$ws.Range("I12").Value = 351.8889
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 351.8889
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -181.8889
$ws.Range("N12").Value = -1340

$ws.Range("H40").Value = 2449.5
$ws.Range("I40").Value = 2449.5
$ws.Range("K40").Value = 2449.5
$ws.Range("M40").Value = -2274.5

$ws.Range("H92").Value = 625.17645
$ws.Range("I92").Value = 623.2143
$ws.Range("K92").Value = 623.2143
$ws.Range("M92").Value = 624.7857

$ws.Range("H125").Value = 533.3333
$ws.Range("I125").Value = 600.3333
$ws.Range("J125").Value = 466.33334
$ws.Range("K125").Value = 5402.9997
$ws.Range("L125").Value = 4197.00006
$ws.Range("M125").Value = -2942.9997
$ws.Range("N125").Value = -9117.00006

$ws.Range("H137").Value = 18183368
$ws.Range("I137").Value = 25001556
$ws.Range("J137").Value = 1533.3334
$ws.Range("K137").Value = 75004668
$ws.Range("L137").Value = 4600.0002
$ws.Range("M137").Value = -75002118
$ws.Range("N137").Value = -9700.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 503.23077
$ws.Range("I5").Value = 52.555557
$ws.Range("K5").Value = 52.555557
$ws.Range("M5").Value = 59.444443

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = $null

$ws.Range("H32").Value = 713639.8
$ws.Range("I32").Value = 855919.1
$ws.Range("K32").Value = 855919.1
$ws.Range("M32").Value = -855632.1

$ws.Range("H45").Value = 1600.6154
$ws.Range("I45").Value = 1588.25
$ws.Range("K45").Value = 1588.25
$ws.Range("M45").Value = -1211.25

$ws.Range("H57").Value = 9664.333000000001
$ws.Range("I57").Value = 9664.333000000001
$ws.Range("K57").Value = 9664.333000000001
$ws.Range("M57").Value = -9180.333000000001

$ws.Range("H61").Value = 5404177.5
$ws.Range("I61").Value = 1896049.1
$ws.Range("J61").Value = 47501720
$ws.Range("K61").Value = 1896049.1
$ws.Range("L61").Value = 47501720
$ws.Range("M61").Value = -1895837.1
$ws.Range("N61").Value = -47502144

$ws.Range("H74").Value = 1406396.6
$ws.Range("I74").Value = 1675266.1
$ws.Range("J74").Value = 8275.4
$ws.Range("K74").Value = 1675266.1
$ws.Range("L74").Value = 8275.4
$ws.Range("M74").Value = -1674392.1
$ws.Range("N74").Value = -10023.4

$ws.Range("H77").Value = 1406396.6
$ws.Range("I77").Value = 1675266.1
$ws.Range("J77").Value = 8275.4
$ws.Range("K77").Value = 8376330.5
$ws.Range("L77").Value = 41377
$ws.Range("M77").Value = -8371962.5
$ws.Range("N77").Value = -50113

$ws.Range("H110").Value = 1262.8334
$ws.Range("I110").Value = 996
$ws.Range("J110").Value = 1529.6666
$ws.Range("K110").Value = 996
$ws.Range("L110").Value = 1529.6666
$ws.Range("M110").Value = 1049
$ws.Range("N110").Value = -5619.6666

$ws.Range("H132").Value = 3561.7368
$ws.Range("I132").Value = 1935.1666
$ws.Range("J132").Value = 6350.143
$ws.Range("K132").Value = 5805.4998
$ws.Range("L132").Value = 19050.429
$ws.Range("M132").Value = -3275.4998
$ws.Range("N132").Value = -24110.429

$ws.Range("H136").Value = 5404177.5
$ws.Range("I136").Value = 1896049.1
$ws.Range("J136").Value = 47501720
$ws.Range("K136").Value = 5688147.300000001
$ws.Range("L136").Value = 142505160
$ws.Range("M136").Value = -5685597.300000001
$ws.Range("N136").Value = -142510260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 503.23077
$ws.Range("I4").Value = 52.555557
$ws.Range("K4").Value = 52.555557
$ws.Range("M4").Value = 62.444443

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = $null

$ws.Range("H113").Value = 7996
$ws.Range("I113").Value = 7996
$ws.Range("K113").Value = 7996
$ws.Range("M113").Value = -5826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 569358.3
$ws.Range("I31").Value = 1376633.8
$ws.Range("J31").Value = 6142.9067
$ws.Range("K31").Value = 1376633.8
$ws.Range("L31").Value = 6142.9067
$ws.Range("M31").Value = -1376338.8
$ws.Range("N31").Value = -6732.9067

$ws.Range("H34").Value = 569358.3
$ws.Range("I34").Value = 1376633.8
$ws.Range("J34").Value = 6142.9067
$ws.Range("K34").Value = 1376633.8
$ws.Range("L34").Value = 6142.9067
$ws.Range("M34").Value = -1376431.8
$ws.Range("N34").Value = -6546.9067

$ws.Range("H107").Value = 701.7778
$ws.Range("I107").Value = 805.7
$ws.Range("J107").Value = 571.875
$ws.Range("K107").Value = 805.7
$ws.Range("L107").Value = 571.875
$ws.Range("M107").Value = 1114.3
$ws.Range("N107").Value = -4411.875

$ws.Range("H132").Value = 2641.147
$ws.Range("I132").Value = 2500.8696
$ws.Range("J132").Value = 2934.4546
$ws.Range("K132").Value = 7502.6088
$ws.Range("L132").Value = 8803.363799999999
$ws.Range("M132").Value = -4972.6088
$ws.Range("N132").Value = -13863.3638

$ws.Range("H134").Value = 3938.7104
$ws.Range("I134").Value = 1627.8235
$ws.Range("J134").Value = 5809.4287
$ws.Range("K134").Value = 4883.470499999999
$ws.Range("L134").Value = 17428.2861
$ws.Range("M134").Value = -2348.470499999999
$ws.Range("N134").Value = -22498.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4664
$ws.Range("I107").Value = 431.625
$ws.Range("K107").Value = 1294.875
$ws.Range("M107").Value = 625.125

$ws.Range("H134").Value = 6803.8823
$ws.Range("I134").Value = 2424.182
$ws.Range("K134").Value = 7272.545999999999
$ws.Range("M134").Value = -2202.545999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 52985.57
$ws.Range("I70").Value = 58483.168
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 58483.168
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = -58213.168
$ws.Range("N70").Value = -20540

$ws.Range("H73").Value = 52985.57
$ws.Range("I73").Value = 58483.168
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 58483.168
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = -57547.168
$ws.Range("N73").Value = -21872

$ws.Range("I80").Value = 7379.2
$ws.Range("J80").Value = 7929.5
$ws.Range("K80").Value = 7379.2
$ws.Range("L80").Value = 7929.5
$ws.Range("M80").Value = -6381.2
$ws.Range("N80").Value = -9925.5

$ws.Range("H83").Value = 7536.4287
$ws.Range("I83").Value = 7379.2
$ws.Range("J83").Value = 7929.5
$ws.Range("K83").Value = 36896
$ws.Range("L83").Value = 39647.5
$ws.Range("M83").Value = -31904
$ws.Range("N83").Value = -49631.5

$ws.Range("H97").Value = 1389.04
$ws.Range("I97").Value = 1211.5454
$ws.Range("J97").Value = 1528.5
$ws.Range("K97").Value = 1211.5454
$ws.Range("L97").Value = 1528.5
$ws.Range("M97").Value = -715.5454
$ws.Range("N97").Value = -2520.5

$ws.Range("H107").Value = 803.9286
$ws.Range("I107").Value = 719.125
$ws.Range("J107").Value = 917
$ws.Range("K107").Value = 719.125
$ws.Range("L107").Value = 917
$ws.Range("M107").Value = 1200.875
$ws.Range("N107").Value = -4757

$ws.Range("H126").Value = 2447.8572
$ws.Range("I126").Value = 2448.7222
$ws.Range("J126").Value = 2442.6667
$ws.Range("K126").Value = 7346.1666
$ws.Range("L126").Value = 7328.000100000001
$ws.Range("M126").Value = -4876.1666
$ws.Range("N126").Value = -12268.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2419.2
$ws.Range("I16").Value = 2419.2
$ws.Range("K16").Value = 2419.2
$ws.Range("M16").Value = -2249.2

$ws.Range("H46").Value = 3614.45
$ws.Range("I46").Value = 2299.6667
$ws.Range("J46").Value = 3846.4707
$ws.Range("K46").Value = 2299.6667
$ws.Range("L46").Value = 3846.4707
$ws.Range("M46").Value = -2111.6667
$ws.Range("N46").Value = -4222.4707

$ws.Range("H82").Value = 4314.1665
$ws.Range("I82").Value = 4475
$ws.Range("J82").Value = 3992.5
$ws.Range("K82").Value = 4475
$ws.Range("L82").Value = 3992.5
$ws.Range("M82").Value = -4114
$ws.Range("N82").Value = -4714.5

$ws.Range("H85").Value = 4314.1665
$ws.Range("I85").Value = 4475
$ws.Range("J85").Value = 3992.5
$ws.Range("K85").Value = 4475
$ws.Range("L85").Value = 3992.5
$ws.Range("M85").Value = -3227
$ws.Range("N85").Value = -6488.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1599.6666
$ws.Range("I107").Value = 2800
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 8400
$ws.Range("L107").Value = 2998.5
$ws.Range("M107").Value = -6480
$ws.Range("N107").Value = -6838.5

$ws.Range("H113").Value = 2898.818
$ws.Range("I113").Value = 2922.75
$ws.Range("J113").Value = 2835
$ws.Range("K113").Value = 8768.25
$ws.Range("L113").Value = 8505
$ws.Range("M113").Value = -6598.25
$ws.Range("N113").Value = -12845

$ws.Range("H136").Value = 1821067.8
$ws.Range("I136").Value = 851791.94
$ws.Range("J136").Value = 5213533
$ws.Range("K136").Value = 2555375.82
$ws.Range("L136").Value = 15640599
$ws.Range("M136").Value = -2552825.82
$ws.Range("N136").Value = -15645699

